$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying TPM data was recomputed; the two rows whose Target cluster
# is "Resolving-Mac" (rows 6 and 11) no longer exist in the new output, so
# delete them (higher row index first so the lower one's index stays valid).
$ws.Rows(11).Delete()
$ws.Rows(6).Delete()

# Refresh the numeric NATMI columns (E:T) for the remaining rows (2-9) with
# the values recomputed from the new TPM data.
# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002491666666666667
$ws.Range("H2").Value = 0.007475
$ws.Range("I2").Value = 0.1635273785303319
$ws.Range("J2").Value = 0.1635273785303318
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1266143333333333
$ws.Range("N2").Value = 0.379843
$ws.Range("O2").Value = 0.0955140869844225
$ws.Range("P2").Value = 0.0955140869844225
$ws.Range("Q2").Value = 0.0003154807138888889
$ws.Range("R2").Value = 0.002839326425
$ws.Range("S2").Value = 0.0156191682572807
$ws.Range("T2").Value = 0.0156191682572807

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002491666666666667
$ws.Range("H3").Value = 0.007475
$ws.Range("I3").Value = 0.1635273785303319
$ws.Range("J3").Value = 0.1635273785303318
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3243313333333334
$ws.Range("N3").Value = 0.972994
$ws.Range("O3").Value = 0.2446659107876707
$ws.Range("P3").Value = 0.2446659107876707
$ws.Range("Q3").Value = 0.0008081255722222223
$ws.Range("R3").Value = 0.00727313015
$ws.Range("S3").Value = 0.04000957500684383
$ws.Range("T3").Value = 0.04000957500684382

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002491666666666667
$ws.Range("H4").Value = 0.007475
$ws.Range("I4").Value = 0.1635273785303319
$ws.Range("J4").Value = 0.1635273785303318
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05644366666666667
$ws.Range("N4").Value = 0.169331
$ws.Range("O4").Value = 0.0425794232437066
$ws.Range("P4").Value = 0.0425794232437066
$ws.Range("Q4").Value = 0.0001406388027777778
$ws.Range("R4").Value = 0.001265749225
$ws.Range("S4").Value = 0.00696290146237682
$ws.Range("T4").Value = 0.00696290146237682

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002491666666666667
$ws.Range("H5").Value = 0.007475
$ws.Range("I5").Value = 0.1635273785303319
$ws.Range("J5").Value = 0.1635273785303318
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8182196666666667
$ws.Range("N5").Value = 2.454659
$ws.Range("O5").Value = 0.6172405789842003
$ws.Range("P5").Value = 0.6172405789842002
$ws.Range("Q5").Value = 0.002038730669444445
$ws.Range("R5").Value = 0.018348576025
$ws.Range("S5").Value = 0.1009357338038305
$ws.Range("T5").Value = 0.1009357338038305

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.01274533333333334
$ws.Range("H6").Value = 0.03823600000000001
$ws.Range("I6").Value = 0.8364726214696682
$ws.Range("J6").Value = 0.8364726214696681
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1266143333333333
$ws.Range("N6").Value = 0.379843
$ws.Range("O6").Value = 0.0955140869844225
$ws.Range("P6").Value = 0.0955140869844225
$ws.Range("Q6").Value = 0.001613741883111111
$ws.Range("R6").Value = 0.014523676948
$ws.Range("S6").Value = 0.0798949187271418
$ws.Range("T6").Value = 0.07989491872714179

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.01274533333333334
$ws.Range("H7").Value = 0.03823600000000001
$ws.Range("I7").Value = 0.8364726214696682
$ws.Range("J7").Value = 0.8364726214696681
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3243313333333334
$ws.Range("N7").Value = 0.972994
$ws.Range("O7").Value = 0.2446659107876707
$ws.Range("P7").Value = 0.2446659107876707
$ws.Range("Q7").Value = 0.004133710953777779
$ws.Range("R7").Value = 0.03720339858400001
$ws.Range("S7").Value = 0.2046563357808269
$ws.Range("T7").Value = 0.2046563357808268

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.01274533333333334
$ws.Range("H8").Value = 0.03823600000000001
$ws.Range("I8").Value = 0.8364726214696682
$ws.Range("J8").Value = 0.8364726214696681
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05644366666666667
$ws.Range("N8").Value = 0.169331
$ws.Range("O8").Value = 0.0425794232437066
$ws.Range("P8").Value = 0.0425794232437066
$ws.Range("Q8").Value = 0.0007193933462222224
$ws.Range("R8").Value = 0.006474540116000001
$ws.Range("S8").Value = 0.03561652178132978
$ws.Range("T8").Value = 0.03561652178132978

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.01274533333333334
$ws.Range("H9").Value = 0.03823600000000001
$ws.Range("I9").Value = 0.8364726214696682
$ws.Range("J9").Value = 0.8364726214696681
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8182196666666667
$ws.Range("N9").Value = 2.454659
$ws.Range("O9").Value = 0.6172405789842003
$ws.Range("P9").Value = 0.6172405789842002
$ws.Range("Q9").Value = 0.01042848239155556
$ws.Range("R9").Value = 0.09385634152400001
$ws.Range("S9").Value = 0.5163048451803698
$ws.Range("T9").Value = 0.5163048451803696
